# Replace the old campaign dates sentence with the new wording,
# everywhere it occurs in the document body.
$d = $word.ActiveDocument

$old = "Datas da campaña de Constelación de Leo 2022: 14-23 de abril, 14-23 de maio"
$new = "Datas da campaña de 2022 que usan Constelación de Leo: 14-23 de abril, 14-23 de maio"

$range = $d.Content
$range.Find.Execute($old, $true, $false, $false, $false, $false, `
                     $true, 1, $false, $new, 2)
